# Update captcha solving functionality for Guarulhos
# Column O ("Observação") results are being refreshed:
#  - Rows previously analysed (all data rows) have their Observação value cleared
#  - Guarulhos rows (102-116) are re-validated and marked as "VÁLIDO"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all previous Observação results for the data rows (2-178)
$ws.Range("O2:O178").Value = ""

# Guarulhos rows (102-116) now validate successfully
$ws.Range("O102:O116").Value = "VÁLIDO"
